# Update the "Förändrad" (Changed) date column (C2:C6) from 2023-10-13 (45212)
# to 2023-10-22 (45221) as described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45221
}
